# Fruta / hortaliza, semanal
# Inserts a new weekly price-report row for "Vega Monumental Concepción - Piña"
# above the existing row 169, shifting all subsequent rows down by one
# (old row 169 becomes 170, ..., old row 249 becomes 250).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 169; Excel shifts rows 169:249 down to 170:250.
$ws.Rows("169:169").Insert()

# Populate the newly inserted row 169 with the new weekly observation.
$ws.Range("A169").Value = 11
$ws.Range("B169").Value = "Vega Monumental Concepción"
$ws.Range("C169").Value = "Bíobío"
$ws.Range("D169").Value = 44992
$ws.Range("E169").Value = 8
$ws.Range("F169").Value = "Fruta"
$ws.Range("G169").Value = 100108
$ws.Range("H169").Value = "Tropicales y subtropicales"
$ws.Range("I169").Value = 100108005
$ws.Range("J169").Value = "Piña"
$ws.Range("K169").Value = "Caramelo"
$ws.Range("L169").Value = "Segunda"
$ws.Range("M169").Value = 300
$ws.Range("N169").Value = 24000
$ws.Range("O169").Value = 25000
$ws.Range("P169").Value = 24500
$ws.Range("Q169").Value = "$/caja 14 unidades"
$ws.Range("R169").Value = "Ecuador"
$ws.Range("S169").Value = 1750
$ws.Range("T169").Value = 14
